$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VEDA_Sets-Proc")

$ws.Range("A21").Value = "IRE"
$ws.Range("B21").Value = "g[_]*"
$ws.Range("F21").Value = "Grid"

$ws.Activate() | Out-Null
$ws.Range("B21").Select() | Out-Null
